$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds the "Date" values for data rows 2-31. The stored value
# "5-7-2012-13" was off by one day (NBA stats source quirk) and needs to
# become the corrected ISO-style date string "2013-05-07" -- stored as
# literal text, not an Excel date serial.
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    if ($cell.Text -eq "5-7-2012-13") {
        $cell.Value = "2013-05-07"
    }
}
